$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B's width is widened to match column A (14.42578125 chars wide).
# The COM ColumnWidth setter snaps to a whole-pixel grid, so 13.66 is the
# input that lands on the closest achievable stored width (14.5).
$ws.Columns.Item(2).ColumnWidth = 13.66

# Update the four rows of paired +/- values in columns A and B.
$ws.Range("A1").Value = 0.037301589504291001
$ws.Range("B1").Value = -0.037301591018868582

$ws.Range("A2").Value = 0.040544987035190647
$ws.Range("B2").Value = -0.040544988547724002

$ws.Range("A3").Value = -0.029367833257236461
$ws.Range("B3").Value = 0.029367831674822539

$ws.Range("A4").Value = 0.023484824666503729
$ws.Range("B4").Value = -0.023484826180156855
